$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solutions")

# Row 28 (occupation group 0 / "All")
$ws.Range("B28:E28").ClearFormats()
$ws.Range("B28").Value = 0.22827110292917499
$ws.Range("C28").Value = 0.36977136334100402
$ws.Range("D28").Value = 0.62817123693064103
$ws.Range("E28").Value = 0.37244128726772702

# Row 29 (occupation group 2)
$ws.Range("B29:E29").ClearFormats()
$ws.Range("B29").Value = 0.22827110292917499
$ws.Range("C29").Value = 0.25326013716395901
$ws.Range("D29").Value = 1.01783384060056
$ws.Range("E29").Value = 0.0348840436770392

# Row 30 (occupation group 1)
$ws.Range("B30:E30").ClearFormats()
$ws.Range("B30").Value = 0.22827110292917499
$ws.Range("C30").Value = 1.0766597280455901
$ws.Range("D30").Value = 0.113785888586552
$ws.Range("E30").Value = 0.0601547142065568

# Recalculate so the dependent formulas (G:I, L:P) refresh their cached values
$excel.Calculate()

# Update the selected cell/range to match the saved view state
$ws.Activate()
$ws.Range("D11").Select()
